$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset the style of A2:A39 to the default (General) style, matching the
# diff where these cells lose their bold/bordered header-ish style (s="1").
$ws.Range("A2:A39").Style = "Normal"

# Add the new "Erlangen" row (row 40) with its view coordinates.
$ws.Cells.Item(40, 1).Value = "Erlangen"
$ws.Cells.Item(40, 2).Value = 49.5928616
$ws.Cells.Item(40, 3).Value = 11.0056
$ws.Cells.Item(40, 4).Value = 10.8556
$ws.Cells.Item(40, 5).Value = 49.4428616
$ws.Cells.Item(40, 6).Value = 11.1556
$ws.Cells.Item(40, 7).Value = 49.7428616
